# Updates the coin Price (D) and Volume(1h) (E) columns, and two pairs of
# rows whose rank order swapped (14/15, 35/36, 49/50), to match the latest
# coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'88.509.17"
$ws.Range("E2").Value = '  +1.70%  '

$ws.Range("D3").Value = "'3.280.28"
$ws.Range("E3").Value = '  -0.87%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = "'215.24"
$ws.Range("E5").Value = '  -1.44%  '

$ws.Range("D6").Value = "'633.84"
$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("D7").Value = "'0.387"
$ws.Range("E7").Value = '  +19.63%  '

$ws.Range("D8").Value = "'0.736"
$ws.Range("E8").Value = '  +20.63%  '

$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("D10").Value = "'3.278.07"
$ws.Range("E10").Value = '  -0.84%  '

$ws.Range("D11").Value = "'0.584"
$ws.Range("E11").Value = '  -2.22%  '

$ws.Range("E12").Value = '  +11.57%  '

$ws.Range("D13").Value = "'0.0000269"
$ws.Range("E13").Value = '  -0.51%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = "'34.79"
$ws.Range("E14").Value = '  +1.14%  '

$ws.Range("B15").Value = 'Toncoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D15").Value = "'5.55"
$ws.Range("E15").Value = '  +3.01%  '

$ws.Range("D16").Value = "'3.881.83"
$ws.Range("E16").Value = '  -0.98%  '

$ws.Range("D17").Value = "'88.230.25"
$ws.Range("E17").Value = '  +1.69%  '

$ws.Range("D18").Value = "'3.286.96"
$ws.Range("E18").Value = '  -0.66%  '

$ws.Range("D19").Value = "'3.18"
$ws.Range("E19").Value = '  +0.71%  '

$ws.Range("D20").Value = "'14.21"
$ws.Range("E20").Value = '  -1.61%  '

$ws.Range("D21").Value = "'442.20"
$ws.Range("E21").Value = '  -2.88%  '

$ws.Range("D22").Value = "'9.03"
$ws.Range("E22").Value = '  +0.12%  '

$ws.Range("D23").Value = "'5.40"
$ws.Range("E23").Value = '  +1.38%  '

$ws.Range("D24").Value = "'7.42"
$ws.Range("E24").Value = '  +0.23%  '

$ws.Range("D25").Value = "'5.32"
$ws.Range("E25").Value = '  -0.69%  '

$ws.Range("D26").Value = "'12.43"
$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").Value = "'0.0000145"
$ws.Range("E27").Value = '  +13.59%  '

$ws.Range("D28").Value = "'3.451.57"
$ws.Range("E28").Value = '  -1.27%  '

$ws.Range("D29").Value = "'77.86"
$ws.Range("E29").Value = '  -0.30%  '

$ws.Range("E30").Value = '  +0.01%  '

$ws.Range("E31").Value = '  -17.42%  '

$ws.Range("D32").Value = "'0.997"
$ws.Range("E32").Value = '  -0.16%  '

$ws.Range("D33").Value = "'8.93"
$ws.Range("E33").Value = '  -3.09%  '

$ws.Range("D34").Value = "'573.51"
$ws.Range("E34").Value = '  -2.86%  '

$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D35").Value = "'7.38"
$ws.Range("E35").Value = '  +12.26%  '

$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").Value = "'1.39"
$ws.Range("E36").Value = '  -8.50%  '

$ws.Range("E37").Value = '  -2.70%  '

$ws.Range("D38").Value = "'0.141"
$ws.Range("E38").Value = '  -6.62%  '

$ws.Range("D39").Value = "'23.12"
$ws.Range("E39").Value = '  -1.02%  '

$ws.Range("D40").Value = "'21.84"
$ws.Range("E40").Value = '  +2.01%  '

$ws.Range("E41").Value = '  -0.07%  '

$ws.Range("D42").Value = "'3.19"
$ws.Range("E42").Value = '  +5.25%  '

$ws.Range("D43").Value = "'0.407"
$ws.Range("E43").Value = '  -2.28%  '

$ws.Range("D44").Value = "'2.05"
$ws.Range("E44").Value = '  +0.35%  '

$ws.Range("E45").Value = '  +0.08%  '

$ws.Range("D46").Value = "'152.50"
$ws.Range("E46").Value = '  -3.82%  '

$ws.Range("E47").Value = '  +21.85%  '

$ws.Range("D48").Value = "'181.21"
$ws.Range("E48").Value = '  -3.76%  '

$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D49").Value = "'1.35"
$ws.Range("E49").Value = '  +0.37%  '

$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").Value = "'44.84"
$ws.Range("E50").Value = '  -4.00%  '

$ws.Range("D51").Value = "'4.28"
$ws.Range("E51").Value = '  +0.60%  '
